$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.060.82"
$ws.Range("E2").Value = "  +1.05%  "

$ws.Range("D3").Value = "2.258.58"
$ws.Range("E3").Value = "  +3.43%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "270.46"
$ws.Range("E5").Value = "  +4.88%  "

$ws.Range("D6").Value = "92.09"
$ws.Range("E6").Value = "  +14.49%  "

$ws.Range("D7").Value = "'0.630"
$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").Value = "0.999"

$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  +6.47%  "

$ws.Range("D10").Value = "45.93"
$ws.Range("E10").Value = "  +7.50%  "

$ws.Range("D11").Value = "0.0974"
$ws.Range("E11").Value = "  +6.50%  "

$ws.Range("D12").Value = "8.35"
$ws.Range("E12").Value = "  +20.60%  "

$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("D14").Value = "2.588.98"
$ws.Range("E14").Value = "  +3.07%  "

$ws.Range("D15").Value = "15.13"
$ws.Range("E15").Value = "  +6.61%  "

$ws.Range("D16").Value = "2.253.24"
$ws.Range("E16").Value = "  +3.86%  "

$ws.Range("D17").Value = "'0.810"
$ws.Range("E17").Value = "  +4.70%  "

$ws.Range("D18").Value = "44.012.35"
$ws.Range("E18").Value = "  +1.07%  "

$ws.Range("D19").Value = "'0.0000105"
$ws.Range("E19").Value = "  +3.19%  "

$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  +3.68%  "

$ws.Range("D21").Value = "'71.00"
$ws.Range("E21").Value = "  +1.75%  "

$ws.Range("E22").Value = "  -1.87%  "

$ws.Range("D23").Value = "235.16"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").Value = "9.18"
$ws.Range("E24").Value = "  +3.37%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").Value = "11.55"
$ws.Range("E26").Value = "  +8.77%  "

$ws.Range("D27").Value = "2.52"
$ws.Range("E27").Value = "  +13.76%  "

$ws.Range("E28").Value = "  +5.50%  "

$ws.Range("D29").Value = "41.35"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").Value = "172.89"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "0.0924"
$ws.Range("E32").Value = "  +6.52%  "

$ws.Range("D33").Value = "21.04"
$ws.Range("E33").Value = "  +3.68%  "

$ws.Range("E34").Value = "  +5.16%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.124"
$ws.Range("E35").Value = "  +2.09%  "

$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  +1.44%  "

$ws.Range("D37").Value = "0.0352"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").Value = "4.31"
$ws.Range("E38").Value = "  -2.80%  "

$ws.Range("E39").Value = "  +25.60%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "13.01"
$ws.Range("E40").Value = "  -0.60%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.229"
$ws.Range("E41").Value = "  +15.70%  "

$ws.Range("D42").Value = "2.21"
$ws.Range("E42").Value = "  +5.78%  "

$ws.Range("D43").Value = "63.83"
$ws.Range("E43").Value = "  +2.16%  "

$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").Value = "0.0998"
$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("D46").Value = "8.39"
$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("D47").Value = "100.43"
$ws.Range("E47").Value = "  +0.00%  "

$ws.Range("E48").Value = "  +5.43%  "

$ws.Range("E49").Value = "  +2.22%  "

$ws.Range("D50").Value = "'0.440"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").Value = "2.479.53"
$ws.Range("E51").Value = "  +3.21%  "
